$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1777.8
$ws.Range("I40").Value = 1550
$ws.Range("J40").Value = 2038.1428
$ws.Range("K40").Value = 1550
$ws.Range("L40").Value = 2038.1428
$ws.Range("M40").Value = -1375
$ws.Range("N40").Value = -2388.1428
$ws.Range("H64").Value = 3034.913
$ws.Range("I64").Value = 2882.2222
$ws.Range("J64").Value = 3133.0715
$ws.Range("K64").Value = 2882.2222
$ws.Range("L64").Value = 3133.0715
$ws.Range("M64").Value = -2634.2222
$ws.Range("N64").Value = -3629.0715
$ws.Range("H67").Value = 3034.913
$ws.Range("I67").Value = 2882.2222
$ws.Range("J67").Value = 3133.0715
$ws.Range("K67").Value = 2882.2222
$ws.Range("L67").Value = 3133.0715
$ws.Range("M67").Value = -2024.2222
$ws.Range("N67").Value = -4849.0715
$ws.Range("H82").Value = 1337.3334
$ws.Range("I82").Value = 1337.3334
$ws.Range("K82").Value = 4012.0002
$ws.Range("M82").Value = -3606.0002
$ws.Range("H85").Value = 1337.3334
$ws.Range("I85").Value = 1337.3334
$ws.Range("K85").Value = 4012.0002
$ws.Range("M85").Value = -2608.0002
$ws.Range("H112").Value = 1465.7894
$ws.Range("J112").Value = 1500
$ws.Range("L112").Value = 4500
$ws.Range("N112").Value = -6716
$ws.Range("H129").Value = 1000.7222
$ws.Range("I129").Value = 502.75
$ws.Range("J129").Value = 1143
$ws.Range("K129").Value = 1508.25
$ws.Range("L129").Value = 3429
$ws.Range("M129").Value = 3491.75
$ws.Range("N129").Value = -13429
$ws.Range("H132").Value = 3377.5095
$ws.Range("I132").Value = 1523.2821
$ws.Range("K132").Value = 4569.846299999999
$ws.Range("M132").Value = -2039.846299999999
$ws.Range("H137").Value = 1135106.2
$ws.Range("I137").Value = 1237.6904
$ws.Range("J137").Value = 7938318
$ws.Range("K137").Value = 3713.0712
$ws.Range("L137").Value = 23814954
$ws.Range("M137").Value = -1163.0712
$ws.Range("N137").Value = -23820054
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1240.6617
$ws.Range("I61").Value = 1071.5741
$ws.Range("K61").Value = 1071.5741
$ws.Range("M61").Value = -859.5741
$ws.Range("H88").Value = 1800.8334
$ws.Range("I88").Value = 1400.7273
$ws.Range("K88").Value = 1400.7273
$ws.Range("M88").Value = -994.7273
$ws.Range("H91").Value = 1800.8334
$ws.Range("I91").Value = 1400.7273
$ws.Range("K91").Value = 1400.7273
$ws.Range("M91").Value = 3.272699999999986
$ws.Range("H136").Value = 1240.6617
$ws.Range("I136").Value = 1071.5741
$ws.Range("K136").Value = 3214.7223
$ws.Range("M136").Value = -664.7223000000004
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1005623.2
$ws.Range("I134").Value = 1744866.5
$ws.Range("J134").Value = 5470.4707
$ws.Range("K134").Value = 5234599.5
$ws.Range("L134").Value = 16411.4121
$ws.Range("M134").Value = -5232064.5
$ws.Range("N134").Value = -21481.4121
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3424.775
$ws.Range("I58").Value = 3570.027
$ws.Range("J58").Value = 1633.3334
$ws.Range("K58").Value = 3570.027
$ws.Range("L58").Value = 1633.3334
$ws.Range("M58").Value = -3367.027
$ws.Range("N58").Value = -2039.3334
$ws.Range("H62").Value = 2423.2727
$ws.Range("I62").Value = 2233.3333
$ws.Range("J62").Value = 2554.7693
$ws.Range("K62").Value = 2233.3333
$ws.Range("L62").Value = 2554.7693
$ws.Range("M62").Value = -1609.3333
$ws.Range("N62").Value = -3802.7693
$ws.Range("H65").Value = 2423.2727
$ws.Range("I65").Value = 2233.3333
$ws.Range("J65").Value = 2554.7693
$ws.Range("K65").Value = 11166.6665
$ws.Range("L65").Value = 12773.8465
$ws.Range("M65").Value = -8046.666499999999
$ws.Range("N65").Value = -19013.8465
$ws.Range("H132").Value = 675322.7
$ws.Range("I132").Value = 1794.28
$ws.Range("K132").Value = 5382.84
$ws.Range("M132").Value = -2852.84
$ws.Range("H136").Value = 3424.775
$ws.Range("I136").Value = 3570.027
$ws.Range("J136").Value = 1633.3334
$ws.Range("K136").Value = 10710.081
$ws.Range("L136").Value = 4900.0002
$ws.Range("M136").Value = -8160.081
$ws.Range("N136").Value = -10000.0002
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 55782.15
$ws.Range("I122").Value = 73843
$ws.Range("K122").Value = 664587
$ws.Range("M122").Value = -662137
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 871.5
$ws.Range("I22").Value = 795
$ws.Range("J22").Value = 880
$ws.Range("K22").Value = 795
$ws.Range("L22").Value = 880
$ws.Range("M22").Value = -500
$ws.Range("N22").Value = -1470
$ws.Range("H27").Value = 871.5
$ws.Range("I27").Value = 795
$ws.Range("J27").Value = 880
$ws.Range("K27").Value = 795
$ws.Range("L27").Value = 880
$ws.Range("M27").Value = -688
$ws.Range("N27").Value = -1094
$ws.Range("H46").Value = 1025.5
$ws.Range("I46").Value = 1275.2
$ws.Range("J46").Value = 847.1429000000001
$ws.Range("K46").Value = 1275.2
$ws.Range("L46").Value = 847.1429000000001
$ws.Range("M46").Value = -1087.2
$ws.Range("N46").Value = -1223.1429
$ws.Range("H82").Value = 2076.923
$ws.Range("I82").Value = 3625
$ws.Range("J82").Value = 1388.8889
$ws.Range("K82").Value = 3625
$ws.Range("L82").Value = 1388.8889
$ws.Range("M82").Value = -3264
$ws.Range("N82").Value = -2110.8889
$ws.Range("H85").Value = 2076.923
$ws.Range("I85").Value = 3625
$ws.Range("J85").Value = 1388.8889
$ws.Range("K85").Value = 3625
$ws.Range("L85").Value = 1388.8889
$ws.Range("M85").Value = -2377
$ws.Range("N85").Value = -3884.8889
$ws.Range("H122").Value = 3703.3333
$ws.Range("I122").Value = 3492.5334
$ws.Range("J122").Value = 3879
$ws.Range("K122").Value = 10477.6002
$ws.Range("L122").Value = 11637
$ws.Range("M122").Value = -8027.600199999999
$ws.Range("H132").Value = 4149.146
$ws.Range("I132").Value = 3853.6667
$ws.Range("J132").Value = 4799.2
$ws.Range("K132").Value = 11561.0001
$ws.Range("L132").Value = 14397.6
$ws.Range("M132").Value = -9031.000100000001
$ws.Range("N132").Value = -19457.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 4941.625
$ws.Range("I107").Value = 6348.9443
$ws.Range("J107").Value = 719.6667
$ws.Range("K107").Value = 19046.8329
$ws.Range("L107").Value = 2159.0001
$ws.Range("M107").Value = -17126.8329
$ws.Range("N107").Value = -5999.0001
$ws.Range("H113").Value = 659.9524
$ws.Range("J113").Value = 665.0909
$ws.Range("L113").Value = 1995.2727
$ws.Range("N113").Value = -6335.2727
$ws.Range("H122").Value = 68844.94500000001
$ws.Range("I122").Value = 6953.6
$ws.Range("J122").Value = 300937.5
$ws.Range("K122").Value = 20860.8
$ws.Range("L122").Value = 902812.5
$ws.Range("M122").Value = -18410.8
